$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Before: row2=ck1.txt(-2,39)  row3=ck2.txt(0,39)  row4=ck3.txt(2,70) [styled]
# After:  row2=ck1.txt(-2,70)  row3=ck3.txt(2,70) [styled]            D5 (new, underlined)
#
# Drop the "ck2.txt" row (row 3) entirely; this shifts the former row 4
# ("ck3.txt", which already carries the bold/"Normale+font1" look) up into
# row 3, and drops the now-unused "ck2.txt" shared string.
$ws.Rows(3).Delete()

# Row 2 ("ck1.txt"): "# Vin Values" goes from 39 to 70, and picks up the
# same formatting as the ck3.txt row beneath it.
$ws.Cells.Item(2, 3).Value = 70
$ws.Cells.Item(2, 3).Style = $ws.Cells.Item(3, 3).Style

# New empty, underlined cell at D5.
$ws.Cells.Item(5, 4).Font.Underline = $true

# Match the saved selection/active cell.
$ws.Cells.Item(5, 4).Select() | Out-Null
